# The document contains a single large marks table (Table 1). Three score
# cells that were previously left at "0" need to be filled in with the
# candidates' actual marks:
#   - Row 11 (MWITI SHARON KARIANKI, DPTE 169/23), column 11 "R/S"  -> 25
#   - Row 47 (KINOTI CAROLINE G.,   DPTE 927/23), column 12 "ER"   -> 44
#   - Row 48 (KAWIRA J. THIAMPURIA, DPTE 985/23), column 12 "ER"   -> 30

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(11, 11).Range.Text = "25"
$t.Cell(47, 12).Range.Text = "44"
$t.Cell(48, 12).Range.Text = "30"
